$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and volume/1h change (E) columns per latest data pull.
# D-column values are forced to text format so numeric-looking strings (e.g. "49.00",
# "7.120") keep their exact textual representation instead of being parsed as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.074.39"
$ws.Range("E2").Value = "  -3.80%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.644.63"
$ws.Range("E3").Value = "  -3.57%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.68"
$ws.Range("E5").Value = "  -2.87%  "
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3904"
$ws.Range("E7").Value = "  -2.44%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3857"
$ws.Range("E8").Value = "  -4.58%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.347"
$ws.Range("E10").Value = "  -8.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "49.00"
$ws.Range("E11").Value = "  -7.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08459"
$ws.Range("E12").Value = "  -4.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.81"
$ws.Range("E13").Value = "  -8.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.120"
$ws.Range("E14").Value = "  -4.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001282"
$ws.Range("E15").Value = "  -5.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.475"
$ws.Range("E16").Value = "  -6.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.646.29"
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.95"
$ws.Range("E18").Value = "  -1.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06973"
$ws.Range("E19").Value = "  -3.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.70"
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.923"
$ws.Range("E21").Value = "  -5.47%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.62"
$ws.Range("E23").Value = "  -4.94%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.080.53"
$ws.Range("E24").Value = "  -3.76%  "
$ws.Range("E25").Value = "  -3.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.704"
$ws.Range("E26").Value = "  -8.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.44"
$ws.Range("E27").Value = "  -4.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "158.24"
$ws.Range("E28").Value = "  -2.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.706"
$ws.Range("E29").Value = "  +3.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "141.51"
$ws.Range("E30").Value = "  -7.21%  "
$ws.Range("E31").Value = "  -13.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.456"
$ws.Range("E32").Value = "  -8.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.821.34"
$ws.Range("E33").Value = "  -3.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.901"
$ws.Range("E34").Value = "  -4.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08019"
$ws.Range("E35").Value = "  -7.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02911"
$ws.Range("E36").Value = "  -8.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9559"
$ws.Range("E37").Value = "  -8.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2692"
$ws.Range("E38").Value = "  -7.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09197"
$ws.Range("E39").Value = "  -5.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.466"
$ws.Range("E40").Value = "  -1.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.928"
$ws.Range("E41").Value = "  -10.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7594"
$ws.Range("E42").Value = "  -8.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.07"
$ws.Range("E43").Value = "  -6.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.02"
$ws.Range("E44").Value = "  -5.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6900"
$ws.Range("E45").Value = "  -6.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.478"
$ws.Range("E46").Value = "  -8.06%  "
$ws.Range("E47").Value = "  -3.74%  "
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08336"
$ws.Range("E49").Value = "  -9.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.57"
$ws.Range("E50").Value = "  -4.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.259"
$ws.Range("E51").Value = "  -10.44%  "
